$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.547.61"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.222.78"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").Value = "2.555.96"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "2.279.93"
$ws.Range("E15").Value = "  -5.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.26%  "
$ws.Range("D18").Value = "44.390.50"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "0.0₃0907"
$ws.Range("E19").Value = "  -6.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -13.27%  "
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "148.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0751"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.39%  "
$ws.Range("E35").Value = "  -3.16%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.28%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0296"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.22%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.71%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "1.817.59"
$ws.Range("E44").Value = "  +3.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.180"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.90%  "
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "68.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "74.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.37%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "93.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.32%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.63%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "13.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.59%  "
